$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = -12.00333333333333
$ws.Range("P2").Value = 180.6266666666667
$ws.Range("I3").Value = -12.00333333333333
$ws.Range("P3").Value = 180.6266666666667
$ws.Range("C4").Value = 18.12333333333333
$ws.Range("I4").Value = -12
$ws.Range("P4").Value = 180.05
$ws.Range("C5").Value = 18.12
$ws.Range("O5").Value = 9.063333333333333
$ws.Range("P5").Value = 180.0533333333333
$ws.Range("B6").Value = -93.96333333333332
$ws.Range("I6").Value = -12.00333333333333
$ws.Range("N6").Value = -46.98333333333333
$ws.Range("O6").Value = 9.060000000000002
$ws.Range("P6").Value = 182.8066666666666
$ws.Range("B7").Value = -93.46333333333332
$ws.Range("C7").Value = 18.12333333333333
$ws.Range("I7").Value = -12
$ws.Range("N7").Value = -46.73
$ws.Range("O7").Value = 9.06
$ws.Range("P7").Value = 319.51
$ws.Range("B8").Value = -93.46333333333332
$ws.Range("C8").Value = 18.12333333333333
$ws.Range("I8").Value = -12
$ws.Range("N8").Value = -46.73
$ws.Range("O8").Value = 9.059999999999997
$ws.Range("P8").Value = 319.51
$ws.Range("B9").Value = -93.96333333333332
$ws.Range("I9").Value = -12.00333333333333
$ws.Range("N9").Value = -46.98333333333333
$ws.Range("O9").Value = 9.060000000000002
$ws.Range("P9").Value = 182.8066666666666
$ws.Range("C10").Value = 18.12
$ws.Range("O10").Value = 9.063333333333333
$ws.Range("P10").Value = 180.0533333333333
$ws.Range("C11").Value = 18.12333333333333
$ws.Range("I11").Value = -12
$ws.Range("P11").Value = 180.05
